$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for 01504b35-... row (G4)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-03 08:49:10"

# zh-cn sheet: Correspond Handoff Datetime (H4) and Correspond Handback DateTime (K4)
# for the 01504b35-... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-09-03 08:49:02"
$wsZhCn.Range("K4").Value = "2016-09-03 08:49:43"

# de-de sheet: Correspond Handback DateTime (K4) for the 01504b35-... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-09-03 08:49:50"
